# Updates crypto price/volume data per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D='30.426.00'; E='  +0.88%  '}
    @{Row=3; D='1.869.42'; E='  +0.45%  '}
    @{Row=4; D='1.001'; DForceText=$true; E='  +0.13%  '}
    @{Row=5; D='246.59'; DForceText=$true; E='  +1.95%  '}
    @{Row=6; D='1.001'; DForceText=$true; E='  +0.12%  '}
    @{Row=7; D='0.4737'; DForceText=$true; E='  +0.80%  '}
    @{Row=8; E='  +1.81%  '}
    @{Row=9; D='0.06497'; DForceText=$true; E='  +0.46%  '}
    @{Row=10; D='22.04'; DForceText=$true; E='  +6.46%  '}
    @{Row=11; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='97.95'; DForceText=$true; E='  +4.40%  '}
    @{Row=12; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.07720'; DForceText=$true; E='  +0.61%  '}
    @{Row=13; D='0.7391'; DForceText=$true; E='  +8.61%  '}
    @{Row=14; D='1.872.19'; E='  +0.61%  '}
    @{Row=15; D='5.115'; DForceText=$true; E='  +1.21%  '}
    @{Row=16; D='273.37'; DForceText=$true; E='  +1.72%  '}
    @{Row=17; D='30.412.07'; E='  +0.86%  '}
    @{Row=18; D='13.39'; DForceText=$true; E='  +0.76%  '}
    @{Row=19; D='0.000007561'; DForceText=$true; E='  +0.19%  '}
    @{Row=20; E='  +0.11%  '}
    @{Row=21; D='2.119.31'; E='  +0.84%  '}
    @{Row=22; D='1.001'; DForceText=$true; E='  +0.15%  '}
    @{Row=23; D='5.238'; DForceText=$true; E='  +1.52%  '}
    @{Row=24; D='6.168'; DForceText=$true; E='  +1.38%  '}
    @{Row=25; D='9.272'; DForceText=$true; E='  -0.53%  '}
    @{Row=26; D='163.84'; DForceText=$true; E='  -1.37%  '}
    @{Row=28; D='1.931'; DForceText=$true; E='  +2.69%  '}
    @{Row=29; D='0.1004'; DForceText=$true; E='  +1.88%  '}
    @{Row=30; E='  -0.82%  '}
    @{Row=31; D='1.508'; DForceText=$true; E='  +0.04%  '}
    @{Row=32; D='4.298'; DForceText=$true; E='  +1.98%  '}
    @{Row=33; D='4.144'; DForceText=$true; E='  +3.74%  '}
    @{Row=34; D='0.04826'; DForceText=$true; E='  +2.79%  '}
    @{Row=35; D='1.119'; DForceText=$true; E='  +0.88%  '}
    @{Row=36; D='0.6966'; DForceText=$true; E='  +1.78%  '}
    @{Row=37; D='2.713'; DForceText=$true; E='  +0.20%  '}
    @{Row=38; E='  +1.78%  '}
    @{Row=39; D='2.753'; DForceText=$true; E='  +1.15%  '}
    @{Row=40; D='6.294'; DForceText=$true; E='  -1.04%  '}
    @{Row=41; D='1.973'; DForceText=$true; E='  +4.80%  '}
    @{Row=42; D='71.56'; DForceText=$true; E='  +2.43%  '}
    @{Row=43; D='0.4182'; DForceText=$true; E='  +3.20%  '}
    @{Row=44; D='1.000'; DForceText=$true; E='  +0.11%  '}
    @{Row=45; D='0.8359'; DForceText=$true; E='  +0.30%  '}
    @{Row=46; D='102.75'; DForceText=$true; E='  +0.83%  '}
    @{Row=47; D='9.301'; DForceText=$true; E='  +0.92%  '}
    @{Row=48; D='7.018'; DForceText=$true; E='  +1.40%  '}
    @{Row=49; E='  +3.32%  '}
    @{Row=50; D='917.46'; DForceText=$true; E='  -0.81%  '}
    @{Row=51; B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.05635'; DForceText=$true; E='  +1.42%  '}
)

foreach ($u in $updates) {
    if ($u.ContainsKey('B')) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey('D')) {
        if ($u.ContainsKey('DForceText')) { $ws.Cells.Item($u.Row, 4).NumberFormat = '@' }
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($u.ContainsKey('E')) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
}
